# Insert a new data row at row 53 (pushing existing rows 53-115 down to 54-116)
# and populate it with the new "Haba" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 53..115 down by one row, creating a blank row 53.
$ws.Rows.Item(53).Insert()

# Fill in the new row 53 with the new record's data.
$ws.Cells.Item(53, 1).Value  = 4
$ws.Cells.Item(53, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(53, 3).Value  = "Los Lagos"
$ws.Cells.Item(53, 4).Value  = 44874
$ws.Cells.Item(53, 5).Value  = 10
$ws.Cells.Item(53, 6).Value  = 100112026
$ws.Cells.Item(53, 7).Value  = "Haba"
$ws.Cells.Item(53, 8).Value  = "Sin especificar"
$ws.Cells.Item(53, 9).Value  = "Primera"
$ws.Cells.Item(53, 10).Value = 70
$ws.Cells.Item(53, 11).Value = 12000
$ws.Cells.Item(53, 12).Value = 12000
$ws.Cells.Item(53, 13).Value = 12000
$ws.Cells.Item(53, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(53, 15).Value = "Región del Maule"
$ws.Cells.Item(53, 16).Value = 480
$ws.Cells.Item(53, 17).Value = 25
$ws.Cells.Item(53, 18).Value = "Hortaliza"
